# The "ppm" column was originally written into column F by mistake; it
# actually belongs in column D (which was left blank). Shift sample_size,
# t_results and significance one column to the left (F<-G, G<-H, H<-I) and
# drop the now-unused column I, then fix up the header labels to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 211

# Grab the old sample_size / t_results / significance block (G:I) in one
# shot and drop it one column to the left (F:H) — preserves numeric vs.
# boolean types per-cell since it's a straight Range.Value copy.
$src = $ws.Range("G2:I$lastRow")
$shifted = $src.Value()
$dest = $ws.Range("F2:H$lastRow")
$dest.Value = $shifted

# The old column I (now duplicated into H) is no longer needed.
$ws.Columns.Item(9).Delete()

# Fix the header row: "ppm" moves to D1, and the remaining labels shift
# left to line up with their shifted data columns.
$ws.Range("D1").Value = "ppm"
$ws.Range("F1").Value = "sample_size"
$ws.Range("G1").Value = "t_results"
$ws.Range("H1").Value = "significance"
